$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

Replace-Text "M2DocEvaluator.caseQuery(M2DocEvaluator.java:516)" "M2DocEvaluator.caseQuery(M2DocEvaluator.java:540)"
Replace-Text "TemplateSwitch.doSwitch(TemplateSwitch.java:172)" "TemplateSwitch.doSwitch(TemplateSwitch.java:186)"
Replace-Text "M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1158)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)"
Replace-Text "TemplateSwitch.doSwitch(TemplateSwitch.java:183)" "TemplateSwitch.doSwitch(TemplateSwitch.java:199)"
Replace-Text "M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:311)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:201)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:266)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)"
Replace-Text "TemplateSwitch.doSwitch(TemplateSwitch.java:246)" "TemplateSwitch.doSwitch(TemplateSwitch.java:279)"
Replace-Text "M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:255)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:705)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:458)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:365)`n`tat sun.reflect.GeneratedMethodAccessor76.invoke(Unknown Source)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)`n`tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)"
